$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.729797666666667
$ws.Range("H2").Value = 5.189393000000001
$ws.Range("I2").Value = 0.06436583050179444
$ws.Range("J2").Value = 0.06436583050179444
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08378199999999998
$ws.Range("N2").Value = 0.251346
$ws.Range("O2").Value = 0.007571394704126512
$ws.Range("P2").Value = 0.007571394704126512
$ws.Range("Q2").Value = 0.1449259081086667
$ws.Range("R2").Value = 1.304333172978
$ws.Range("S2").Value = 0.0004873391081879911
$ws.Range("T2").Value = 0.0004873391081879911
$ws.Range("G3").Value = 1.729797666666667
$ws.Range("H3").Value = 5.189393000000001
$ws.Range("I3").Value = 0.06436583050179444
$ws.Range("J3").Value = 0.06436583050179444
$ws.Range("M3").Value = 0.07352966666666667
$ws.Range("O3").Value = 0.006644889460697858
$ws.Range("P3").Value = 0.006644889460697857
$ws.Range("Q3").Value = 0.1271914458307778
$ws.Range("R3").Value = 1.144723012477
$ws.Range("S3").Value = 0.0004277038287304386
$ws.Range("T3").Value = 0.0004277038287304385
$ws.Range("G4").Value = 1.729797666666667
$ws.Range("H4").Value = 5.189393000000001
$ws.Range("I4").Value = 0.06436583050179444
$ws.Range("J4").Value = 0.06436583050179444
$ws.Range("O4").Value = 0.9857837158351757
$ws.Range("P4").Value = 0.9857837158351755
$ws.Range("Q4").Value = 18.869125940335
$ws.Range("R4").Value = 169.822133463015
$ws.Range("S4").Value = 0.06345078756487602
$ws.Range("T4").Value = 0.063450787564876
$ws.Range("I5").Value = 0.2200595722726403
$ws.Range("J5").Value = 0.2200595722726403
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08378199999999998
$ws.Range("N5").Value = 0.251346
$ws.Range("O5").Value = 0.007571394704126512
$ws.Range("P5").Value = 0.007571394704126512
$ws.Range("Q5").Value = 0.4954854633426666
$ws.Range("R5").Value = 4.459369170083999
$ws.Range("S5").Value = 0.001666157880097414
$ws.Range("T5").Value = 0.001666157880097414
$ws.Range("I6").Value = 0.2200595722726403
$ws.Range("J6").Value = 0.2200595722726403
$ws.Range("M6").Value = 0.07352966666666667
$ws.Range("O6").Value = 0.006644889460697858
$ws.Range("P6").Value = 0.006644889460697857
$ws.Range("S6").Value = 0.001462271532520146
$ws.Range("T6").Value = 0.001462271532520146
$ws.Range("I7").Value = 0.2200595722726403
$ws.Range("J7").Value = 0.2200595722726403
$ws.Range("O7").Value = 0.9857837158351757
$ws.Range("P7").Value = 0.9857837158351755
$ws.Range("S7").Value = 0.2169311428600227
$ws.Range("T7").Value = 0.2169311428600227
$ws.Range("H8").Value = 57.69206699999999
$ws.Range("I8").Value = 0.7155745972255653
$ws.Range("J8").Value = 0.7155745972255653
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.08378199999999998
$ws.Range("N8").Value = 0.251346
$ws.Range("O8").Value = 0.007571394704126512
$ws.Range("P8").Value = 0.007571394704126512
$ws.Range("Q8").Value = 1.611185585798
$ws.Range("R8").Value = 14.500670272182
$ws.Range("S8").Value = 0.005417897715841107
$ws.Range("T8").Value = 0.005417897715841107
$ws.Range("H9").Value = 57.69206699999999
$ws.Range("I9").Value = 0.7155745972255653
$ws.Range("J9").Value = 0.7155745972255653
$ws.Range("M9").Value = 0.07352966666666667
$ws.Range("O9").Value = 0.006644889460697858
$ws.Range("P9").Value = 0.006644889460697857
$ws.Range("S9").Value = 0.004754914099447273
$ws.Range("T9").Value = 0.004754914099447272
$ws.Range("H10").Value = 57.69206699999999
$ws.Range("I10").Value = 0.7155745972255653
$ws.Range("J10").Value = 0.7155745972255653
$ws.Range("O10").Value = 0.9857837158351757
$ws.Range("P10").Value = 0.9857837158351755
$ws.Range("S10").Value = 0.7054017854102769
$ws.Range("T10").Value = 0.7054017854102769

Write-Output "Applied 88 updated TPM values to Fgf1-Fgfr4 sheet"
